$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Diciembre")

$ws.Range("B4").Value = 5163
$ws.Range("E4").Value = 2467

$ws.Range("B9").Value = 14493
$ws.Range("C9").Value = 677

$ws.Range("B12").Value = 7616
$ws.Range("E12").Value = 1301

$ws.Range("B14").Value = 16360
$ws.Range("C14").Value = 481
